# Update the "FUELS" sheet of the SIN LCA_infrastructure workbook:
# the Natural Gas (NG) row's PEN and CO2 figures are refreshed to match
# ecoinvent 3.4 "market for natural gas, burned in gas motor, for
# storage_GLO_2017_Allocation, cut-off", and the reference column is
# updated to cite that new source string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")

# PEN (C2): was a literal 1.403, now computed from the ecoinvent components
$ws.Range("C2").Formula = "=1.1767+0.0019487+0.0000015726"

# CO2 (D2): updated figure
$ws.Range("D2").Value = 0.06682

# reference (F2): point at the new ecoinvent source string
$ws.Range("F2").Value = "ecoinvent 3.4 - market for natural gas, burned in gas motor, for storage_GLO_2017_Allocation, cut-off"

# move the sheet's cursor/selection to C2, matching the saved view state
$ws.Range("C2").Select()
